# Applies the "correct svm pic and table, and finish 'Different Libraries of
# SVM'" edit: updates two chart titles to call out LibSVM, tightens the
# value-axis scaling on the three SVM comparison charts (slides 19, 20, 22)
# so the bars are easier to read, and leaves everything else untouched.

$p = $ppt.ActivePresentation

# --- Slide 19 : "Letter Recognition" SVM chart -> "Letter Recognition (LibSVM)"
$c19 = $p.Slides.Item(19).Shapes.Item(2).Chart
$c19.ChartTitle.TextFrame.TextRange.Text = "Letter Recognition (LibSVM)"
$ax19 = $c19.Axes(2)
$ax19.MinimumScale = 50.0
$ax19.MajorUnit = 10.0

# --- Slide 20 : "Letter Recognition (SMO vs LibSVM)" chart - title unchanged,
# only the value-axis scaling is tightened to match.
$c20 = $p.Slides.Item(20).Shapes.Item(2).Chart
$ax20 = $c20.Axes(2)
$ax20.MinimumScale = 50.0
$ax20.MajorUnit = 10.0

# --- Slide 22 : "Digit Classification" SVM chart -> "Digit Classification (LibSVM)"
$c22 = $p.Slides.Item(22).Shapes.Item(2).Chart
$c22.ChartTitle.TextFrame.TextRange.Text = "Digit Classification (LibSVM)"
$ax22 = $c22.Axes(2)
$ax22.MinimumScale = 50.0
$ax22.MajorUnit = 10.0
